# Updated cryptos list values (Price column D, Volume(1h) column E).
# Price cells that look like plain numbers get a leading apostrophe so Excel
# keeps them as text (matching the original text-formatted cells) instead of
# auto-converting them to numeric values and dropping formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.670.26'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Value = '1.802.27'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").Value = '''231.14'
$ws.Range("E5").Value = '  -2.39%  '
$ws.Range("D6").Value = '''0.5939'
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").Value = '''0.2777'
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").Value = '''0.06842'
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").Value = '''23.37'
$ws.Range("E10").Value = '  -2.21%  '
$ws.Range("D11").Value = '''0.07533'
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("D12").Value = '1.807.21'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").Value = '''4.712'
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("D14").Value = '''0.6259'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '2.047.42'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '''0.000009231'
$ws.Range("E16").Value = '  -8.72%  '
$ws.Range("D17").Value = '''75.32'
$ws.Range("E17").Value = '  -5.22%  '
$ws.Range("D18").Value = '28.630.43'
$ws.Range("D19").Value = '''5.467'
$ws.Range("E19").Value = '  -7.62%  '
$ws.Range("D20").Value = '''1.004'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = '''210.84'
$ws.Range("E21").Value = '  -7.66%  '
$ws.Range("D22").Value = '''11.44'
$ws.Range("E22").Value = '  -3.30%  '
$ws.Range("D23").Value = '''6.847'
$ws.Range("E23").Value = '  -2.71%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '''154.43'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '''7.836'
$ws.Range("D27").Value = '''0.1274'
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("D28").Value = '''16.39'
$ws.Range("E28").Value = '  -1.50%  '
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").Value = '''0.06264'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("D32").Value = '''3.767'
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("D33").Value = '''3.738'
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("E35").Value = '  -6.72%  '
$ws.Range("D36").Value = '''0.6398'
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").Value = '''2.502'
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").Value = '''2.722'
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").Value = '''0.01710'
$ws.Range("E39").Value = '  -2.38%  '
$ws.Range("D40").Value = '''6.393'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").Value = '1.137.98'
$ws.Range("E41").Value = '  -6.47%  '
$ws.Range("D42").Value = '''0.8651'
$ws.Range("E42").Value = '  -7.17%  '
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = '''100.58'
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '1.964.38'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '''60.45'
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("E47").Value = '  -6.01%  '
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("D49").Value = '''8.334'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = '''0.05470'
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("E51").Value = '  -1.32%  '
